# Apply "transfered citations to the excel sheet" edit:
#  - Insert a new worksheet "COM.others" right after "Disc.Typical.Landscape"
#  - Populate it with the transferred citation text
#  - Update the "Disc.Typical.Landscape" sheet's selected cell

$wb = $excel.ActiveWorkbook

# Update the selection on "Disc.Typical.Landscape" before it stops being the
# active sheet, so the new cursor position (A27) is recorded for that sheet.
$landscapeSheet = $wb.Worksheets.Item("Disc.Typical.Landscape")
$landscapeSheet.Activate()
$landscapeSheet.Range("A27").Select()

# Insert the new "COM.others" worksheet right after "Disc.Typical.Landscape".
$afterSheet = $wb.Worksheets.Item("Disc.Typical.Landscape")
$comSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$comSheet.Name = "COM.others"

# Transfer the citation text into column A of the new sheet.
$comSheet.Range("A2").Value = "other COM, / nucleus membrane citations (Harper, Golubovskaya, & Cande, 2004; Naranjo & Corredor, 2008; Scherthan et al., 1996; Zickler & Kleckner, 2016)."
$comSheet.Range("A4").Value = "telomere guided movements  citations: (Bass et al.,"
$comSheet.Range("A5").Value = "2000; Chacon, Delivani, & Tolic, 2016; Curtis, Lukaszewski, & Chrzastek, 1991; Ding, Yamamoto, Haraguchi, & Hiraoka, 2004; Gerton &"
$comSheet.Range("A6").Value = "Hawley, 2005; Lee, Conrad, & Dresser, 2012; Lefrancois, Rockmill,"
$comSheet.Range("A7").Value = "Xie, Roeder, & Snyder, 2016; Page & Hawley, 2003)"
$comSheet.Range("A9").Value = "Repair of DSB - CO"
$comSheet.Range("A10").Value = "(Anderson & Stack, 2005; Bass et al., 2000; Brown et al., 2005; Croft & Jones, 1989; Higgins, Osman, Jones, & Franklin, 2014; Klutstein & Cooper, 2014, Lukaszewski, 1997; Pratto et al., 2014; Viera, Santos, & Rufas, 2009;"
$comSheet.Range("A11").Value = "Xiang, Miller, Ross, Alvarado, & Hawley, 2014)."

# Leave the cursor on F8 of the new sheet, matching the source workbook.
$comSheet.Range("F8").Select()
